$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AccountSet": rows 2-6 get cyclically rotated (row2->6, row3->2stays..)
# Concretely or2=Interest, or3=Principal, or4=PrevStmt, or5=CurrStmt, or6=Checking
# becomes: nr2=Checking, nr3=CurrStmt, nr4=PrevStmt(same row, date becomes text),
#          nr5=Principal(date becomes text), nr6=Interest
# Capture all the old row values first so the rewrite doesn't clobber data we
# still need to read.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AccountSet")

$or2_A = $ws.Cells.Item(2, 1).Value()
$or2_B = $ws.Cells.Item(2, 2).Value()
$or2_C = $ws.Cells.Item(2, 3).Value()
$or2_D = $ws.Cells.Item(2, 4).Value()
$or2_E = $ws.Cells.Item(2, 5).Value()

$or3_A = $ws.Cells.Item(3, 1).Value()
$or3_B = $ws.Cells.Item(3, 2).Value()
$or3_C = $ws.Cells.Item(3, 3).Value()
$or3_D = $ws.Cells.Item(3, 4).Value()
$or3_E = $ws.Cells.Item(3, 5).Value()
$or3_F = $ws.Cells.Item(3, 6).Value()
$or3_G = $ws.Cells.Item(3, 7).Value()
$or3_H = $ws.Cells.Item(3, 8).Value()
$or3_I = $ws.Cells.Item(3, 9).Value()
$or3_J = $ws.Cells.Item(3, 10).Value()

$or5_A = $ws.Cells.Item(5, 1).Value()
$or5_B = $ws.Cells.Item(5, 2).Value()
$or5_C = $ws.Cells.Item(5, 3).Value()
$or5_D = $ws.Cells.Item(5, 4).Value()
$or5_E = $ws.Cells.Item(5, 5).Value()

$or6_A = $ws.Cells.Item(6, 1).Value()
$or6_B = $ws.Cells.Item(6, 2).Value()
$or6_C = $ws.Cells.Item(6, 3).Value()
$or6_D = $ws.Cells.Item(6, 4).Value()
$or6_E = $ws.Cells.Item(6, 5).Value()

# New row 2 <- old row 6 (Checking)
$ws.Cells.Item(2, 1).Value = $or6_A
$ws.Cells.Item(2, 2).Value = $or6_B
$ws.Cells.Item(2, 3).Value = $or6_C
$ws.Cells.Item(2, 4).Value = $or6_D
$ws.Cells.Item(2, 5).Value = $or6_E

# New row 3 <- old row 5 (Credit: Curr Stmt Bal) ; F3:J3 stay blank
$ws.Cells.Item(3, 1).Value = $or5_A
$ws.Cells.Item(3, 2).Value = $or5_B
$ws.Cells.Item(3, 3).Value = $or5_C
$ws.Cells.Item(3, 4).Value = $or5_D
$ws.Cells.Item(3, 5).Value = $or5_E
$ws.Range("F3:J3").ClearContents()

# Row 4 stays "Credit: Prev Stmt Bal" but Billing_Start_Dt (F4) becomes text
$ws.Cells.Item(4, 6).Value = "'20000102"
$ws.Cells.Item(4, 6).ClearFormats()

# New row 5 <- old row 3 (test loan: Principal Balance), F5 becomes text date
$ws.Cells.Item(5, 1).Value = $or3_A
$ws.Cells.Item(5, 2).Value = $or3_B
$ws.Cells.Item(5, 3).Value = $or3_C
$ws.Cells.Item(5, 4).Value = $or3_D
$ws.Cells.Item(5, 5).Value = $or3_E
$ws.Cells.Item(5, 6).Value = "'20000102"
$ws.Cells.Item(5, 6).ClearFormats()
$ws.Cells.Item(5, 7).Value = $or3_G
$ws.Cells.Item(5, 8).Value = $or3_H
$ws.Cells.Item(5, 9).Value = $or3_I
$ws.Cells.Item(5, 10).Value = $or3_J

# New row 6 <- old row 2 (test loan: Interest)
$ws.Cells.Item(6, 1).Value = $or2_A
$ws.Cells.Item(6, 2).Value = $or2_B
$ws.Cells.Item(6, 3).Value = $or2_C
$ws.Cells.Item(6, 4).Value = $or2_D
$ws.Cells.Item(6, 5).Value = $or2_E

# ---------------------------------------------------------------------------
# Sheet "BudgetSet": Start_Date/End_Date (cols A,B) on rows 2-4 switch from
# numeric yyyymmdd values to literal text of the same digits.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("BudgetSet")
foreach ($r in 2..4) {
    foreach ($c in 1..2) {
        $cell = $ws2.Cells.Item($r, $c)
        $digits = [string]([int]$cell.Value())
        $cell.Value = "'" + $digits
        $cell.ClearFormats()
    }
}

# ---------------------------------------------------------------------------
# Sheet "config": Start_Date_YYYYMMDD/End_Date_YYYYMMDD (A2,B2) switch from
# numeric to literal text as well.
# ---------------------------------------------------------------------------
$ws9 = $wb.Worksheets.Item("config")
$cellA2 = $ws9.Cells.Item(2, 1)
$digitsA2 = [string]([int]$cellA2.Value())
$cellA2.Value = "'" + $digitsA2
$cellA2.ClearFormats()

$cellB2 = $ws9.Cells.Item(2, 2)
$digitsB2 = [string]([int]$cellB2.Value())
$cellB2.Value = "'" + $digitsB2
$cellB2.ClearFormats()
